$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain-looking number but must stay text
# (matches original inlineStr cell type) - force Text number format first.
$ws.Range("D2").Value = "60.013.94"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").Value = "3.195.70"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.97"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.64"
$ws.Range("E6").Value = "  +4.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("E10").Value = "  +1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.434"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").Value = "3.745.88"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.137"
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.77"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").Value = "60.021.59"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").Value = "3.208.29"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.26"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.20"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "370.30"
$ws.Range("E21").Value = "  -1.63%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.523"
$ws.Range("E23").Value = "  -1.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.64"
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("E26").Value = "  +3.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").Value = "0.0₃0878"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.47"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.12"
$ws.Range("E31").Value = "  +1.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.28"
$ws.Range("E32").Value = "  +2.15%  "
$ws.Range("E33").Value = "  +4.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.08"
$ws.Range("E35").Value = "  +1.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.37"
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.53"
$ws.Range("E37").Value = "  +5.70%  "
$ws.Range("D38").Value = "2.781.85"
$ws.Range("E38").Value = "  +4.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0315"
$ws.Range("E39").Value = "  +9.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0710"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.718"
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("D46").Value = "3.236.02"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.982"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("E50").Value = "  +5.18%  "
